{"js": "// The duplicated-bookmark error paragraphs get a \"    <---\" marker prefixed\n// to the existing red/bold warning text, so they stand out from the\n// (correctly) bookmarked example just above them. The run formatting\n// (bold + red) must be preserved; only the text content of the run changes.\n\nconst startResults = context.document.body.search(\"Can't start duplicated bookmark bookmark1\", { matchCase: true });\nstartResults.load(\"items\");\nawait context.sync();\n\nif (startResults.items.length > 0) {\n  startResults.items[0].insertText(\n    \"    <---Can't start duplicated bookmark bookmark1\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\nconst endResults = context.document.body.search(\"Can't end already closed bookmark bookmark1\", { matchCase: true });\nendResults.load(\"items\");\nawait context.sync();\n\nif (endResults.items.length > 0) {\n  endResults.items[0].insertText(\n    \"    <---Can't end already closed bookmark bookmark1\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# The duplicated-bookmark error paragraphs get a \"    <---\" marker prefixed\n# to the existing red/bold warning text, so they stand out from the\n# (correctly) bookmarked example just above them. The run formatting\n# (bold + red) must be preserved; only the text content of the run changes.\n\n$d = $word.ActiveDocument\n\n$range1 = $d.Content\n$find1 = $range1.Find\n$find1.ClearFormatting()\n$find1.Text = \"Can't start duplicated bookmark bookmark1\"\n$found1 = $find1.Execute()\nif ($found1) {\n    $range1.Text = \"    <---Can't start duplicated bookmark bookmark1\"\n}\n\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.ClearFormatting()\n$find2.Text = \"Can't end already closed bookmark bookmark1\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $range2.Text = \"    <---Can't end already closed bookmark bookmark1\"\n}\n"}
